$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5285.3335
$ws.Range("I19").Value = 8163.231
$ws.Range("J19").Value = 608.75
$ws.Range("K19").Value = 8163.231
$ws.Range("L19").Value = 608.75
$ws.Range("M19").Value = -7988.231
$ws.Range("N19").Value = -958.75

$ws.Range("H33").Value = 679.2308
$ws.Range("I33").Value = 442.6
$ws.Range("J33").Value = 1101.7858
$ws.Range("K33").Value = 442.6
$ws.Range("L33").Value = 1101.7858
$ws.Range("M33").Value = -213.6
$ws.Range("N33").Value = -1559.7858

$ws.Range("H38").Value = 678.2105
$ws.Range("J38").Value = 2300
$ws.Range("L38").Value = 6900
$ws.Range("N38").Value = -7644

$ws.Range("H41").Value = 9804.909
$ws.Range("I41").Value = 481.625
$ws.Range("J41").Value = 34667
$ws.Range("K41").Value = 481.625
$ws.Range("L41").Value = 34667
$ws.Range("M41").Value = -41.625
$ws.Range("N41").Value = -35547

$ws.Range("H76").Value = 8267.179
$ws.Range("I76").Value = 12623.417
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 12623.417
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -12308.417
$ws.Range("N76").Value = -5630

$ws.Range("H79").Value = 8267.179
$ws.Range("I79").Value = 12623.417
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 12623.417
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -11531.417
$ws.Range("N79").Value = -7184

$ws.Range("H98").Value = 1181.8077
$ws.Range("I98").Value = 1046
$ws.Range("J98").Value = 1487.375
$ws.Range("K98").Value = 1046
$ws.Range("L98").Value = 1487.375
$ws.Range("M98").Value = 452
$ws.Range("N98").Value = -4483.375

$ws.Range("H113").Value = 2926.8157
$ws.Range("I113").Value = 2556.3845
$ws.Range("K113").Value = 2556.3845
$ws.Range("M113").Value = 697.6154999999999

$ws.Range("H122").Value = 1181.8077
$ws.Range("I122").Value = 1046
$ws.Range("J122").Value = 1487.375
$ws.Range("K122").Value = 3138
$ws.Range("L122").Value = 4462.125
$ws.Range("M122").Value = -688
$ws.Range("N122").Value = -9362.125

$ws.Range("H132").Value = 3884.3809
$ws.Range("I132").Value = 1842.8889
$ws.Range("J132").Value = 16133.333
$ws.Range("K132").Value = 5528.6667
$ws.Range("L132").Value = 48399.999
$ws.Range("M132").Value = -2998.6667
$ws.Range("N132").Value = -53459.999

$ws.Range("H137").Value = 36319.766
$ws.Range("I137").Value = 56247.42
$ws.Range("J137").Value = 1899.2727
$ws.Range("K137").Value = 168742.26
$ws.Range("L137").Value = 5697.8181
$ws.Range("M137").Value = -166192.26
$ws.Range("N137").Value = -10797.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1362.4615
$ws.Range("I45").Value = 1220.3334
$ws.Range("J45").Value = 1484.2858
$ws.Range("K45").Value = 1220.3334
$ws.Range("L45").Value = 1484.2858
$ws.Range("M45").Value = -843.3334
$ws.Range("N45").Value = -2238.2858

$ws.Range("H132").Value = 3201.2307
$ws.Range("I132").Value = 2987.353
$ws.Range("J132").Value = 3605.2222
$ws.Range("K132").Value = 8962.059000000001
$ws.Range("L132").Value = 10815.6666
$ws.Range("M132").Value = -6432.059000000001
$ws.Range("N132").Value = -15875.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = ""
$ws.Range("N75").Value = ""

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = ""
$ws.Range("N78").Value = ""

$ws.Range("H82").Value = 25823.334
$ws.Range("I82").Value = 8333.333000000001
$ws.Range("J82").Value = 34568.332
$ws.Range("K82").Value = 8333.333000000001
$ws.Range("L82").Value = 34568.332
$ws.Range("M82").Value = -7950.333000000001
$ws.Range("N82").Value = -35334.332

$ws.Range("H85").Value = 25823.334
$ws.Range("I85").Value = 8333.333000000001
$ws.Range("J85").Value = 34568.332
$ws.Range("K85").Value = 8333.333000000001
$ws.Range("L85").Value = 34568.332
$ws.Range("M85").Value = -7007.333000000001
$ws.Range("N85").Value = -37220.332

$ws.Range("H112").Value = 30400
$ws.Range("J112").Value = 30400
$ws.Range("L112").Value = 30400
$ws.Range("N112").Value = -33354

$ws.Range("H134").Value = 22230.371
$ws.Range("I134").Value = 26780.244
$ws.Range("J134").Value = 7880.769
$ws.Range("K134").Value = 80340.73199999999
$ws.Range("L134").Value = 23642.307
$ws.Range("M134").Value = -77805.73199999999
$ws.Range("N134").Value = -28712.307

$ws.Range("H139").Value = 44695
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 44695
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 44695
$ws.Range("M139").Value = ""
$ws.Range("N139").Value = -54975

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 63.875
$ws.Range("I7").Value = 46.27778
$ws.Range("J7").Value = 116.666664
$ws.Range("K7").Value = 46.27778
$ws.Range("L7").Value = 116.666664
$ws.Range("M7").Value = 66.72221999999999
$ws.Range("N7").Value = -342.666664

$ws.Range("H99").Value = 61786.35
$ws.Range("I99").Value = 113388.445
$ws.Range("J99").Value = 3734
$ws.Range("K99").Value = 113388.445
$ws.Range("L99").Value = 3734
$ws.Range("M99").Value = -111890.445
$ws.Range("N99").Value = -6730

$ws.Range("H126").Value = 61786.35
$ws.Range("I126").Value = 113388.445
$ws.Range("J126").Value = 3734
$ws.Range("K126").Value = 340165.335
$ws.Range("L126").Value = 11202
$ws.Range("M126").Value = -337695.335
$ws.Range("N126").Value = -16142

$ws.Range("H134").Value = 2044.1428
$ws.Range("I134").Value = 1208.7142
$ws.Range("J134").Value = 2879.5715
$ws.Range("K134").Value = 3626.1426
$ws.Range("L134").Value = 8638.7145
$ws.Range("M134").Value = -1091.1426
$ws.Range("N134").Value = -13708.7145

$ws.Range("H141").Value = 24114.084
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 24114.084
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 24114.084
$ws.Range("M141").Value = ""
$ws.Range("N141").Value = -34474.084

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 135.23529
$ws.Range("J12").Value = 195.72728
$ws.Range("L12").Value = 587.18184
$ws.Range("N12").Value = -933.18184

$ws.Range("H43").Value = 3200
$ws.Range("J43").Value = 3200
$ws.Range("L43").Value = 9600
$ws.Range("N43").Value = -9828

$ws.Range("H58").Value = 2631.111
$ws.Range("I58").Value = 800
$ws.Range("J58").Value = 2860
$ws.Range("K58").Value = 2400
$ws.Range("L58").Value = 8580
$ws.Range("M58").Value = -2272
$ws.Range("N58").Value = -8836

$ws.Range("H131").Value = 854.2241
$ws.Range("I131").Value = 523.6
$ws.Range("J131").Value = 923.1042
$ws.Range("K131").Value = 1570.8
$ws.Range("L131").Value = 2769.3126
$ws.Range("M131").Value = 3469.2
$ws.Range("N131").Value = -12849.3126

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 22850.166
$ws.Range("J139").Value = 22850.166
$ws.Range("L139").Value = 22850.166
$ws.Range("N139").Value = -33130.166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2164.6667
$ws.Range("I7").Value = 1971.0834
$ws.Range("J7").Value = 2422.7778
$ws.Range("K7").Value = 1971.0834
$ws.Range("L7").Value = 2422.7778
$ws.Range("M7").Value = -1859.0834
$ws.Range("N7").Value = -2646.7778

$ws.Range("H16").Value = 2650.0715
$ws.Range("I16").Value = 2675.0833
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 2675.0833
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -2505.0833
$ws.Range("N16").Value = -2840

$ws.Range("H41").Value = 9220
$ws.Range("I41").Value = 8500
$ws.Range("J41").Value = 9700
$ws.Range("K41").Value = 8500
$ws.Range("L41").Value = 9700
$ws.Range("M41").Value = -8062
$ws.Range("N41").Value = -10576

$ws.Range("H126").Value = 2164.6667
$ws.Range("I126").Value = 1971.0834
$ws.Range("J126").Value = 2422.7778
$ws.Range("K126").Value = 5913.2502
$ws.Range("L126").Value = 7268.3334
$ws.Range("M126").Value = -3443.2502
$ws.Range("N126").Value = -12208.3334

$ws.Range("H139").Value = 41225
$ws.Range("J139").Value = 41225
$ws.Range("L139").Value = 41225
$ws.Range("N139").Value = -51505

$ws.Range("H140").Value = 46718.445
$ws.Range("J140").Value = 46718.445
$ws.Range("L140").Value = 46718.445
$ws.Range("N140").Value = -57078.445

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 421
$ws.Range("I107").Value = 430.72726
$ws.Range("J107").Value = 385.33334
$ws.Range("K107").Value = 1292.18178
$ws.Range("L107").Value = 1156.00002
$ws.Range("M107").Value = 627.8182200000001
$ws.Range("N107").Value = -4996.000019999999

$ws.Range("H109").Value = 30300
$ws.Range("J109").Value = 30300
$ws.Range("L109").Value = 30300
$ws.Range("N109").Value = -33074

$ws.Range("H132").Value = 19238.465
$ws.Range("I132").Value = 28960.555
$ws.Range("J132").Value = 3329.5908
$ws.Range("K132").Value = 86881.66500000001
$ws.Range("L132").Value = 9988.7724
$ws.Range("M132").Value = -84351.66500000001
$ws.Range("N132").Value = -15048.7724

$ws.Range("H136").Value = 19609476
$ws.Range("I136").Value = 37038304
$ws.Range("J136").Value = 2046.0416
$ws.Range("K136").Value = 111114912
$ws.Range("L136").Value = 6138.1248
$ws.Range("M136").Value = -111112362
$ws.Range("N136").Value = -11238.1248
